# Fixed some little bugs in NetLiquidity script.
#
# RRPONTSYD.xlsx gets five new daily observations appended to the "Data"
# sheet, and the "SeriesInfo" sheet's metadata is refreshed to match the
# newer FRED pull.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) "Data" sheet: append rows 444-448 (2023-07-13 .. 2023-07-19).
# ---------------------------------------------------------------------
$wsData = $wb.Worksheets.Item("Data")

# Clone the formatting (date number format / font / border / alignment)
# of the last existing data row onto the new rows before filling values.
$wsData.Range("A443").Copy()
$wsData.Range("A444:A448").PasteSpecial(-4122)

$wsData.Range("A444").Value = 45120
$wsData.Range("B444").Value = 1767.432

$wsData.Range("A445").Value = 45121
$wsData.Range("B445").Value = 1740.777

$wsData.Range("A446").Value = 45124
$wsData.Range("B446").Value = 1728.322

$wsData.Range("A447").Value = 45125
$wsData.Range("B447").Value = 1716.862

$wsData.Range("A448").Value = 45126
$wsData.Range("B448").Value = 1732.804

# ---------------------------------------------------------------------
# 2) "SeriesInfo" sheet: refresh the realtime/observation/update stamps.
# ---------------------------------------------------------------------
$wsInfo = $wb.Worksheets.Item("SeriesInfo")

# These replacement strings ("2023-07-20", "2023-07-19", ...) look like
# dates, so a plain `.Value =` assignment would get auto-converted into a
# date serial. Force the cell to text first, write the value, then strip
# the formatting back off so the cell ends up as a plain, unstyled string
# -- same as the other metadata cells in this column.
function Set-TextValue($range, [string]$text) {
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.ClearFormats()
}

Set-TextValue $wsInfo.Range("B3") "2023-07-20"
Set-TextValue $wsInfo.Range("B4") "2023-07-20"
Set-TextValue $wsInfo.Range("B7") "2023-07-19"
Set-TextValue $wsInfo.Range("B14") "2023-07-19 13:01:03-05"
